$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 182, pushing existing rows 182-196 down to 183-197.
$ws.Rows.Item(182).Insert()

# Populate the new row 182 with the new record (Región de Coquimbo).
$ws.Cells.Item(182, 1).Value = 10
$ws.Cells.Item(182, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(182, 3).Value = "La Araucanía"
$ws.Cells.Item(182, 4).Value = 44826
$ws.Cells.Item(182, 5).Value = 9
$ws.Cells.Item(182, 6).Value = 100112012
$ws.Cells.Item(182, 7).Value = "Espinaca"
$ws.Cells.Item(182, 8).Value = "Sin especificar"
$ws.Cells.Item(182, 9).Value = "Primera"
$ws.Cells.Item(182, 10).Value = 50
$ws.Cells.Item(182, 11).Value = 10000
$ws.Cells.Item(182, 12).Value = 10000
$ws.Cells.Item(182, 13).Value = 10000
$ws.Cells.Item(182, 14).Value = "`$/docena de atados"
$ws.Cells.Item(182, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(182, 16).Value = 3333
$ws.Cells.Item(182, 17).Value = 3
$ws.Cells.Item(182, 18).Value = "Hortaliza"
